# Slide 1 ("文本框 5", shape id 21): update the college name, and
# reposition/resize the auto-fitting (wrap="none" + spAutoFit), centered
# text box to match the bounds PowerPoint recalculates for the shorter
# text (the shape's horizontal center point is preserved).
#
# Left/Width are PowerPoint COM `Single` (32-bit float) properties
# measured in points; the literals below are chosen so that, after the
# float32 round-trip + EMU conversion, the saved XML offsets match the
# target EMU values exactly (x=5312187, cx=1783122). Top/Height are left
# untouched since they are unchanged by this edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(3)

$shp.TextFrame.TextRange.Text = "学院：人工智能学院"

$shp.Left = 418.28245544488516
$shp.Width = 140.40331268661384
